$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (raw OOXML width target in comment) ---
$colsTo8 = @(2,3,7,9,10,11,12,13,15,16,17,22,24,26,27,28,29,30,34)
foreach ($col in $colsTo8) {
  $ws.Columns.Item($col).ColumnWidth = 7.17
}
$ws.Columns.Item(20).ColumnWidth = 8.17

# --- Replace data rows 2-5 with new values ---
$arr = New-Object 'object[,]' 4,34
$arr[0,0] = 45133.50694444445
$arr[0,1] = 13.798
$arr[0,2] = 9.137
$arr[0,3] = 3.527
$arr[0,4] = 29.879
$arr[0,5] = 22.444
$arr[0,6] = 10.657
$arr[0,7] = 31.967
$arr[0,8] = 16.87
$arr[0,9] = 6.742
$arr[0,10] = 10.011
$arr[0,11] = 11.733
$arr[0,12] = 12.516
$arr[0,13] = 3.497
$arr[0,14] = 10.903
$arr[0,15] = 14.966
$arr[0,16] = 9.704000000000001
$arr[0,17] = 3.059
$arr[0,18] = 1.672
$arr[0,19] = 158.575
$arr[0,20] = 30.182
$arr[0,21] = 10.064
$arr[0,22] = 19.331
$arr[0,23] = 9.888999999999999
$arr[0,24] = 2.87
$arr[0,25] = 17.102
$arr[0,26] = 8.888999999999999
$arr[0,27] = 8.15
$arr[0,28] = 9.673999999999999
$arr[0,29] = 12.058
$arr[0,30] = 3.072
$arr[0,31] = 28.934
$arr[0,32] = 5.407
$arr[0,33] = 12.581
$arr[1,0] = 45133.51388888889
$arr[1,1] = 22.949
$arr[1,2] = 16.738
$arr[1,3] = 1.98
$arr[1,4] = 50.14
$arr[1,5] = 40.24
$arr[1,6] = 17.938
$arr[1,7] = 67.355
$arr[1,8] = 27.923
$arr[1,9] = 12.189
$arr[1,10] = 17.951
$arr[1,11] = 20.031
$arr[1,12] = 21.301
$arr[1,13] = 5.797
$arr[1,14] = 18.046
$arr[1,15] = 25.513
$arr[1,16] = 15.45
$arr[1,17] = 1.506
$arr[1,18] = 1.215
$arr[1,19] = 267.368
$arr[1,20] = 50.439
$arr[1,21] = 16.657
$arr[1,22] = 33.575
$arr[1,23] = 17.551
$arr[1,24] = 2.987
$arr[1,25] = 33.521
$arr[1,26] = 14.713
$arr[1,27] = 13.16
$arr[1,28] = 15.5
$arr[1,29] = 20.882
$arr[1,30] = 1.198
$arr[1,31] = 61.623
$arr[1,32] = 9.279
$arr[1,33] = 20.825
$arr[2,0] = 45133.52083333334
$arr[2,1] = 21.524
$arr[2,2] = 15.88
$arr[2,3] = 1.475
$arr[2,4] = 47.041
$arr[2,5] = 38.053
$arr[2,6] = 16.855
$arr[2,7] = 66.61
$arr[2,8] = 26.178
$arr[2,9] = 11.54
$arr[2,10] = 17.004
$arr[2,11] = 18.831
$arr[2,12] = 20.022
$arr[2,13] = 5.435
$arr[2,14] = 16.918
$arr[2,15] = 24.007
$arr[2,16] = 14.384
$arr[2,17] = 1.039
$arr[2,18] = 0.965
$arr[2,19] = 250.207
$arr[2,20] = 47.339
$arr[2,21] = 15.616
$arr[2,22] = 31.662
$arr[2,23] = 16.577
$arr[2,24] = 2.629
$arr[2,25] = 32.357
$arr[2,26] = 13.794
$arr[2,27] = 12.285
$arr[2,28] = 14.454
$arr[2,29] = 19.688
$arr[2,30] = 0.746
$arr[2,31] = 60.737
$arr[2,32] = 8.741
$arr[2,33] = 19.524
$arr[3,0] = 45133.52777777778
$arr[3,1] = 16.25
$arr[3,2] = 12.02
$arr[3,3] = 1.09
$arr[3,4] = 35.55
$arr[3,5] = 28.75
$arr[3,6] = 12.73
$arr[3,7] = 51.67
$arr[3,8] = 19.78
$arr[3,9] = 8.73
$arr[3,10] = 12.83
$arr[3,11] = 14.24
$arr[3,12] = 15.15
$arr[3,13] = 4.11
$arr[3,14] = 12.78
$arr[3,15] = 18.14
$arr[3,16] = 10.88
$arr[3,17] = 0.78
$arr[3,18] = 0.71
$arr[3,19] = 187.25
$arr[3,20] = 35.8
$arr[3,21] = 11.8
$arr[3,22] = 23.94
$arr[3,23] = 12.53
$arr[3,24] = 1.99
$arr[3,25] = 24.83
$arr[3,26] = 10.42
$arr[3,27] = 9.289999999999999
$arr[3,28] = 10.92
$arr[3,29] = 14.89
$arr[3,30] = 0.54
$arr[3,31] = 47.09
$arr[3,32] = 6.6
$arr[3,33] = 14.75

$ws.Range("A2:AH5").Value = $arr

# --- Remove old row 6 (dataset now has 4 data rows instead of 5) ---
$ws.Rows.Item(6).Delete()
